# The author's edit repositions/resizes the "Chart 9" graphic frame on
# slide 11 (the pie-chart graphicFrame whose cNvPr id="10"). Everything
# else in the diff (c16:uniqueId, a16:creationId, p14:modId) is internal,
# non-deterministic bookkeeping metadata that PowerPoint regenerates on
# its own and that is not exposed/settable through the PowerPoint object
# model (VBA/COM) - so we only reproduce the observable geometry change.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)

$chart = $null
foreach ($shp in $s.Shapes) {
    if ($shp.HasChart -and $shp.Name -eq "Chart 9") {
        $chart = $shp
        break
    }
}
if ($chart -eq $null) {
    $chart = $s.Shapes.Item(7)
}

# Target EMUs: off x=1666875 y=1778547, ext cx=6584458 cy=4207915
# (points below are tuned so the float32 round-trip used internally lands
# on those exact EMU values instead of one EMU short).
$chart.Left   = 131.25
$chart.Top    = 140.04307086614173
$chart.Width  = 518.4613037109375
$chart.Height = 331.3319396972656
